# Arrange for Excel & Add Acl Logging
# Adds four new Key/en/ko-KR/ko rows (Source, Destination, Protocol, Length)
# to the "locale" worksheet, immediately after the existing data (row 153),
# and updates the active selection to reflect where the editor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New localization rows (154-157) -------------------------------------
# Columns: A = Key, B = en, C = ko-KR, D = ko

$ws.Range("A154").Value = "Source"
$ws.Range("B154").Value = "Source"
$ws.Range("C154").Value = "출발지"
$ws.Range("D154").Value = "출발지"

$ws.Range("A155").Value = "Destination"
$ws.Range("B155").Value = "Destination"
$ws.Range("C155").Value = "도착지"
$ws.Range("D155").Value = "도착지"

$ws.Range("A156").Value = "Protocol"
$ws.Range("B156").Value = "Protocol"
$ws.Range("C156").Value = "프로토콜"
$ws.Range("D156").Value = "프로토콜"

$ws.Range("A157").Value = "Length"
$ws.Range("B157").Value = "Length"
$ws.Range("C157").Value = "길이"
$ws.Range("D157").Value = "길이"

# --- Restore the author's final cursor position ---------------------------
$ws.Range("C153").Select()
